# Update "Maximum Capacity Factor.xlsx" per latest data refresh ("updated 4.0 files and mdl")

$wb = $excel.ActiveWorkbook

# --- About sheet: bump the "last updated" date stamp (C1) ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45392

# --- MCF sheet: raise several plant-type capacity factors to 1 (100%) ---
$wsMCF = $wb.Worksheets.Item("MCF")

$wsMCF.Range("B2").Value  = 1   # nuclear
$wsMCF.Range("B3").Value  = 1   # natural gas steam turbine
$wsMCF.Range("B4").Value  = 1   # natural gas combined cycle
$wsMCF.Range("B6").Value  = 1   # hydro
$wsMCF.Range("B10").Value = 1   # biomass
$wsMCF.Range("B11").Value = 1   # geothermal
$wsMCF.Range("B12").Value = 1   # petroleum
$wsMCF.Range("B13").Value = 1   # natural gas peaker
$wsMCF.Range("B14").Value = 1   # lignite
$wsMCF.Range("B16").Value = 1   # heavy or residual fuel oil
$wsMCF.Range("B17").Value = 1   # municipal solid waste
$wsMCF.Range("B18").Value = 1   # (row 18)

# Leave the active selection on the MCF sheet at B17, matching the saved file state
$wsMCF.Activate()
$wsMCF.Range("B17").Select()
